$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new entry: link text cell E7 with hyperlink to the FirstBadVersion.cs source
$ws.Range("E7").Value = "https://github.com/Gershon-Tadepalli/DS-Algo/blob/master/DS-AlgoPractice/DS-AlgoLibrary/SortingAndSearching/FirstBadVersion.cs"

$ws.Hyperlinks.Add(
    $ws.Range("E7"),
    "https://github.com/Gershon-Tadepalli/DS-Algo/blob/master/DS-AlgoPractice/DS-AlgoLibrary/SortingAndSearching/FirstBadVersion.cs"
)

$ws.Range("E7").Style = "Hyperlink"

# Reset scroll position so topLeftCell reverts to A1 and select E7 like original
$ws.Range("A1").Select()
$ws.Range("E7").Select()
